$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B, F, G hold date-like text (e.g. "2024-06-05") that must stay as
# literal text rather than being auto-converted to Excel date serials, so
# force Text format on those columns before writing any values.
$dateCols = @("B", "F", "G")
foreach ($col in $dateCols) {
    $ws.Range("${col}2:${col}24").NumberFormat = "@"
}

# Row 2
$ws.Range("A2").Value = "DB"
$ws.Range("B2").Value = "2024-06-05"
$ws.Range("C2").Value = "디비금융스팩12호"
$ws.Range("D2").Value = "DB"
$ws.Range("E2").Value = "DB"
$ws.Range("F2").Value = "2024-06-11"
$ws.Range("G2").Value = "2024-06-18"
$ws.Range("H2").Value = 10000
$ws.Range("I2").Value = 5000000
$ws.Range("J2").Value = 2000
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 100

# Row 3
$ws.Range("A3").Value = "KB"
$ws.Range("B3").Value = "2024-04-18"
$ws.Range("C3").Value = "제일엠앤에스"
$ws.Range("D3").Value = "KB"
$ws.Range("E3").Value = "KB"
$ws.Range("F3").Value = "2024-04-23"
$ws.Range("G3").Value = "2024-04-30"
$ws.Range("H3").Value = 52800
$ws.Range("I3").Value = 2400000
$ws.Range("J3").Value = 22000
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 100

# Row 4
$ws.Range("A4").Value = "KB"
$ws.Range("B4").Value = "2024-04-23"
$ws.Range("C4").Value = "민테크"
$ws.Range("D4").Value = "KB"
$ws.Range("E4").Value = "KB"
$ws.Range("F4").Value = "2024-04-26"
$ws.Range("G4").Value = "2024-05-03"
$ws.Range("H4").Value = 31500
$ws.Range("I4").Value = 3000000
$ws.Range("J4").Value = 10500
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 100

# Row 5
$ws.Range("A5").Value = "KB"
$ws.Range("B5").Value = "2024-05-07"
$ws.Range("C5").Value = "KB제28호스팩"
$ws.Range("D5").Value = "KB"
$ws.Range("E5").Value = "KB"
$ws.Range("F5").Value = "2024-05-10"
$ws.Range("G5").Value = "2024-05-17"
$ws.Range("H5").Value = 10000
$ws.Range("I5").Value = 5000000
$ws.Range("J5").Value = 2000
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 100

# Row 6
$ws.Range("A6").Value = "KB"
$ws.Range("B6").Value = "2024-04-25"
$ws.Range("C6").Value = "HD현대마린솔루션"
$ws.Range("D6").Value = "KB, 유비에스리미티드(영업소), 제이피모간회사 서울지점"
$ws.Range("E6").Value = "KB, 유비에스리미티드(영업소), 제이피모간회사 서울지점, 신한, 하나, 대신, 삼성"
$ws.Range("F6").Value = "2024-04-30"
$ws.Range("G6").Value = "2024-05-08"
$ws.Range("H6").Value = 215255.4
$ws.Range("I6").Value = 8900000
$ws.Range("J6").Value = 83400
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 29

# Row 7
$ws.Range("A7").Value = "NH"
$ws.Range("B7").Value = "2024-05-07"
$ws.Range("C7").Value = "아이씨티케이"
$ws.Range("D7").Value = "NH"
$ws.Range("E7").Value = "NH"
$ws.Range("F7").Value = "2024-05-10"
$ws.Range("G7").Value = "2024-05-17"
$ws.Range("H7").Value = 39400
$ws.Range("I7").Value = 1970000
$ws.Range("J7").Value = 20000
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 100

# Row 8
$ws.Range("A8").Value = "SK"
$ws.Range("B8").Value = "2024-04-23"
$ws.Range("C8").Value = "SK증권제12호스팩"
$ws.Range("D8").Value = "SK"
$ws.Range("E8").Value = "SK"
$ws.Range("F8").Value = "2024-04-26"
$ws.Range("G8").Value = "2024-05-07"
$ws.Range("H8").Value = 6000
$ws.Range("I8").Value = 3000000
$ws.Range("J8").Value = 2000
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 100

# Row 9
$ws.Range("A9").Value = "대신"
$ws.Range("B9").Value = "2024-06-05"
$ws.Range("C9").Value = "라메디텍"
$ws.Range("D9").Value = "대신"
$ws.Range("E9").Value = "대신"
$ws.Range("F9").Value = "2024-06-11"
$ws.Range("G9").Value = "2024-06-17"
$ws.Range("H9").Value = 20768
$ws.Range("I9").Value = 1298000
$ws.Range("J9").Value = 16000
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 100

# Row 10
$ws.Range("A10").Value = "대신"
$ws.Range("B10").Value = "2024-04-25"
$ws.Range("C10").Value = "HD현대마린솔루션"
$ws.Range("D10").Value = "KB, 유비에스리미티드(영업소), 제이피모간회사 서울지점"
$ws.Range("E10").Value = "KB, 유비에스리미티드(영업소), 제이피모간회사 서울지점, 신한, 하나, 대신, 삼성"
$ws.Range("F10").Value = "2024-04-30"
$ws.Range("G10").Value = "2024-05-08"
$ws.Range("H10").Value = 18556.5
$ws.Range("I10").Value = 8900000
$ws.Range("J10").Value = 83400
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 2.5

# Row 11
$ws.Range("A11").Value = "미래"
$ws.Range("B11").Value = "2024-06-10"
$ws.Range("C11").Value = "미래에셋비전스팩5호"
$ws.Range("D11").Value = "미래"
$ws.Range("E11").Value = "미래"
$ws.Range("F11").Value = "2024-06-13"
$ws.Range("G11").Value = "2024-06-19"
$ws.Range("H11").Value = 9500
$ws.Range("I11").Value = 4750000
$ws.Range("J11").Value = 2000
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 100

# Row 12
$ws.Range("A12").Value = "미래"
$ws.Range("B12").Value = "2024-05-20"
$ws.Range("C12").Value = "미래에셋비전스팩4호"
$ws.Range("D12").Value = "미래"
$ws.Range("E12").Value = "미래"
$ws.Range("F12").Value = "2024-05-23"
$ws.Range("G12").Value = "2024-05-29"
$ws.Range("H12").Value = 13300
$ws.Range("I12").Value = 6650000
$ws.Range("J12").Value = 2000
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 100

# Row 13
$ws.Range("A13").Value = "삼성"
$ws.Range("B13").Value = "2024-06-03"
$ws.Range("C13").Value = "그리드위즈"
$ws.Range("D13").Value = "삼성"
$ws.Range("E13").Value = "삼성"
$ws.Range("F13").Value = "2024-06-07"
$ws.Range("G13").Value = "2024-06-14"
$ws.Range("H13").Value = 56000
$ws.Range("I13").Value = 1400000
$ws.Range("J13").Value = 40000
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 100

# Row 14
$ws.Range("A14").Value = "삼성"
$ws.Range("B14").Value = "2024-05-13"
$ws.Range("C14").Value = "노브랜드"
$ws.Range("D14").Value = "삼성"
$ws.Range("E14").Value = "삼성"
$ws.Range("F14").Value = "2024-05-17"
$ws.Range("G14").Value = "2024-05-23"
$ws.Range("H14").Value = 16800
$ws.Range("I14").Value = 1200000
$ws.Range("J14").Value = 14000
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 100

# Row 15
$ws.Range("A15").Value = "삼성"
$ws.Range("B15").Value = "2024-04-25"
$ws.Range("C15").Value = "HD현대마린솔루션"
$ws.Range("D15").Value = "KB, 유비에스리미티드(영업소), 제이피모간회사 서울지점"
$ws.Range("E15").Value = "KB, 유비에스리미티드(영업소), 제이피모간회사 서울지점, 신한, 하나, 대신, 삼성"
$ws.Range("F15").Value = "2024-04-30"
$ws.Range("G15").Value = "2024-05-08"
$ws.Range("H15").Value = 18556.5
$ws.Range("I15").Value = 8900000
$ws.Range("J15").Value = 83400
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 2.5

# Row 16
$ws.Range("A16").Value = "신한"
$ws.Range("B16").Value = "2024-04-25"
$ws.Range("C16").Value = "HD현대마린솔루션"
$ws.Range("D16").Value = "KB, 유비에스리미티드(영업소), 제이피모간회사 서울지점"
$ws.Range("E16").Value = "KB, 유비에스리미티드(영업소), 제이피모간회사 서울지점, 신한, 하나, 대신, 삼성"
$ws.Range("F16").Value = "2024-04-30"
$ws.Range("G16").Value = "2024-05-08"
$ws.Range("H16").Value = 74226
$ws.Range("I16").Value = 8900000
$ws.Range("J16").Value = 83400
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 10

# Row 17
$ws.Range("A17").Value = "유비에스"
$ws.Range("B17").Value = "2024-04-25"
$ws.Range("C17").Value = "HD현대마린솔루션"
$ws.Range("D17").Value = "KB, 유비에스리미티드(영업소), 제이피모간회사 서울지점"
$ws.Range("E17").Value = "KB, 유비에스리미티드(영업소), 제이피모간회사 서울지점, 신한, 하나, 대신, 삼성"
$ws.Range("F17").Value = "2024-04-30"
$ws.Range("G17").Value = "2024-05-08"
$ws.Range("H17").Value = 170719.8
$ws.Range("I17").Value = 8900000
$ws.Range("J17").Value = 83400
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 23

# Row 18
$ws.Range("A18").Value = "유안타"
$ws.Range("B18").Value = "2024-04-22"
$ws.Range("C18").Value = "유안타제16호스팩"
$ws.Range("D18").Value = "유안타"
$ws.Range("E18").Value = "유안타"
$ws.Range("F18").Value = "2024-04-25"
$ws.Range("G18").Value = "2024-05-02"
$ws.Range("H18").Value = 10300
$ws.Range("I18").Value = 5150000
$ws.Range("J18").Value = 2000
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 100

# Row 19
$ws.Range("A19").Value = "제이피모간회사"
$ws.Range("B19").Value = "2024-04-25"
$ws.Range("C19").Value = "HD현대마린솔루션"
$ws.Range("D19").Value = "KB, 유비에스리미티드(영업소), 제이피모간회사 서울지점"
$ws.Range("E19").Value = "KB, 유비에스리미티드(영업소), 제이피모간회사 서울지점, 신한, 하나, 대신, 삼성"
$ws.Range("F19").Value = "2024-04-30"
$ws.Range("G19").Value = "2024-05-08"
$ws.Range("H19").Value = 170719.8
$ws.Range("I19").Value = 8900000
$ws.Range("J19").Value = 83400
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 23

# Row 20
$ws.Range("A20").Value = "하나"
$ws.Range("B20").Value = "2024-04-25"
$ws.Range("C20").Value = "HD현대마린솔루션"
$ws.Range("D20").Value = "KB, 유비에스리미티드(영업소), 제이피모간회사 서울지점"
$ws.Range("E20").Value = "KB, 유비에스리미티드(영업소), 제이피모간회사 서울지점, 신한, 하나, 대신, 삼성"
$ws.Range("F20").Value = "2024-04-30"
$ws.Range("G20").Value = "2024-05-08"
$ws.Range("H20").Value = 74226
$ws.Range("I20").Value = 8900000
$ws.Range("J20").Value = 83400
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 10

# Row 21
$ws.Range("A21").Value = "한국"
$ws.Range("B21").Value = "2024-04-24"
$ws.Range("C21").Value = "코칩"
$ws.Range("D21").Value = "한국"
$ws.Range("E21").Value = "한국"
$ws.Range("F21").Value = "2024-04-29"
$ws.Range("G21").Value = "2024-05-07"
$ws.Range("H21").Value = 27000
$ws.Range("I21").Value = 1500000
$ws.Range("J21").Value = 18000
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 100

# Row 22
$ws.Range("A22").Value = "한국"
$ws.Range("B22").Value = "2024-06-10"
$ws.Range("C22").Value = "씨어스테크놀로지"
$ws.Range("D22").Value = "한국"
$ws.Range("E22").Value = "한국"
$ws.Range("F22").Value = "2024-06-13"
$ws.Range("G22").Value = "2024-06-19"
$ws.Range("H22").Value = 22100
$ws.Range("I22").Value = 1300000
$ws.Range("J22").Value = 17000
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 100

# Row 23
$ws.Range("A23").Value = "한국"
$ws.Range("B23").Value = "2024-04-22"
$ws.Range("C23").Value = "디앤디파마텍"
$ws.Range("D23").Value = "한국"
$ws.Range("E23").Value = "한국"
$ws.Range("F23").Value = "2024-04-25"
$ws.Range("G23").Value = "2024-05-02"
$ws.Range("H23").Value = 36300
$ws.Range("I23").Value = 1100000
$ws.Range("J23").Value = 33000
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 100

# Row 24
$ws.Range("A24").Value = "한국"
$ws.Range("B24").Value = "2024-06-10"
$ws.Range("C24").Value = "한국제14호스팩"
$ws.Range("D24").Value = "한국"
$ws.Range("E24").Value = "한국"
$ws.Range("F24").Value = "2024-06-13"
$ws.Range("G24").Value = "2024-06-19"
$ws.Range("H24").Value = 8000
$ws.Range("I24").Value = 4000000
$ws.Range("J24").Value = 2000
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 100

# Clean up: drop the per-cell style override added by the Text number format
# above so date cells end up with no explicit style, matching plain data cells.
foreach ($col in $dateCols) {
    $ws.Range("${col}2:${col}24").Style = "Normal"
}